$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency list values (price/volume refresh + row reshuffle)
# Each value is written with a leading apostrophe to force text interpretation
# (matching the source data, which stores numeric-looking values as text),
# then the style is reset to Normal so no extra "quote prefix" formatting is
# left behind on the cell.

$ws.Range("D2").Value = "'29.058.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -4.08%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.963.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -6.25%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.32%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'326.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -4.29%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.25%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4992"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -5.95%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.4212"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -3.90%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'53.81"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.42%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.09149"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.23%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'1.097"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -6.69%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'23.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -6.53%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.986.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.69%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'7.861"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -8.40%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'6.432"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -6.49%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +0.31%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.00001100"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -4.89%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'91.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -9.92%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.06670"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.79%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'19.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -9.24%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.33%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.940"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -6.20%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'29.077.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.99%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'11.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.00%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -1.48%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Value = "'Monero"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'156.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.46%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "'EthereumClassic"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'20.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -5.62%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "'InternetComputer(DFINITY)"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'6.151"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -11.37%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'LidoDAOToken"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'2.266"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -9.59%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "'BitcoinCash"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'126.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -5.13%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "'ImmutableX"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'1.042"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -7.72%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "'Stellar"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'0.09851"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -6.51%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "'ARBITRUM"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'1.534"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -7.99%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "'Filecoin"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'5.773"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -7.58%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'HuobiToken"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'3.681"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -5.90%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'VeChain"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'0.02425"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -7.63%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'TrustWalletToken"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'1.298"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -3.20%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'FraxShare"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'8.923"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -11.26%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'Hedera"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'0.06303"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -6.84%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'TheSandbox"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'0.6457"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -7.17%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'Aptos"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'11.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -9.12%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'Algorand"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.1991"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -10.13%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'Frax"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'1.004"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.28%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'Decentraland"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.6209"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -8.44%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'EnergySwap"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'13.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -5.51%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'NEARProtocol"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'2.172"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -7.62%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'WEMIXTOKEN"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'1.284"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.36%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'PancakeSwap"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'3.463"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.72%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'BabyDogeCoin"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.00000000333"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -4.33%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'Cronos"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.06903"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -5.10%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'ThetaToken"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'1.102"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -8.68%  "
$ws.Range("E51").Style = "Normal"
